$wb = $excel.ActiveWorkbook

# Update the "Correspond Handoff Datetime" / "Correspond Handback DateTime"
# values in the per-language handback status sheets, simulating a fresh
# report generation run (new timestamps a short while after the old ones).

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E4").Value = "2016-03-20 08:35:15"
$wsZhCn.Range("H4").Value = "2016-03-20 08:35:36"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E4").Value = "2016-03-20 08:35:18"
$wsDeDe.Range("H4").Value = "2016-03-20 08:35:42"
